# correccion transaccion movimientos eprepago
# Updates the "Eprepago" sheet test-data: the card/document number used in
# rows 2 and 3 changes from 93221453 to 22452521, and the card suffix on
# row 2 changes from *9344 to *7848. Also refreshes the saved selection on
# the sheet (A2 instead of K12, no frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eprepago")

# --- Correct the test data values -----------------------------------------
# Use a leading apostrophe so these numeric-looking values are kept as text
# (matching the existing "quote prefix" text formatting already on these
# cells), rather than being reinterpreted as numbers.
$ws.Range("B2").Value = "'22452521"
$ws.Range("B3").Value = "'22452521"
$ws.Range("N2").Value = "*7848"

# --- Update the sheet's saved view/selection --------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
